$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add two new "current" model runs for 2035 and two for 2050 (NoProject 06b
# and Plan 07b), demoting the previous "current" rows (05b / 06b) back to
# plain history rows. We insert new rows right after each existing "current"
# row, then fill them in; the previously-current row's Status (I) and Alias
# (V) cells are cleared so only the newest run is marked current.
#
# Original row numbers (before any insert) of the 4 "current" rows that are
# being demoted:
#   110 = 2035 NoProject 05b
#   117 = 2035 Plan 06b
#   124 = 2050 NoProject 05b
#   131 = 2050 Plan 06b
#
# We insert from the bottom up so that the original row numbers for the
# insertion points further up the sheet stay valid while we work.
# ---------------------------------------------------------------------------

# --- 4) 2050 Plan 07b -- insert after original row 131 (2050 Plan 06b) ---
$ws.Rows(132).Insert()
$ws.Range("A132").Value2 = 2050
$ws.Range("B132").Value2 = "2050_TM160_DBP_Plan_07b"
$ws.Range("C132").Value2 = "RTP2025"
$ws.Range("D132").Value2 = "DBP"
$ws.Range("F132").Value2 = "Update landuse/popsyn (fix too many persons); tweak bike back"
$ws.Range("G132").Value2 = "M:\urban_modeling\baus\PBA50Plus\PBA50Plus_DraftBlueprint\PBA50Plus_Draft_Blueprint_v6"
$ws.Range("H132").Value2 = "PBA50Plus_Draft_Blueprint_v6"
$ws.Range("I132").Value2 = "current"
$ws.Range("J132").Value2 = "yes"
$ws.Range("K132").Value2 = "BlueprintNetworks_v18\net_2050_Blueprint"
$ws.Range("L132").Value2 = "model2-c"
$ws.Range("M132").Value2 = "https://app.asana.com/0/1182463234225195/1207127280475846/f"
$ws.Range("N132").Value2 = 19.57
$ws.Range("O132").Value2 = "na"
$ws.Range("P132").Value2 = "na"
$ws.Range("Q132").Value2 = 0.87
$ws.Range("R132").Value2 = 0.78
$ws.Range("S132").Value2 = 83.3
$ws.Range("T132").Value2 = 0
$ws.Range("U132").Value2 = 72
$ws.Range("V132").Value2 = "2050 Plan"

# Demote the previous 2050 Plan "current" row (now/still row 131)
$ws.Range("I131").Clear()
$ws.Range("V131").Clear()

# --- 3) 2050 NoProject 06b -- insert after original row 124 (2050 NoProject 05b) ---
$ws.Rows(125).Insert()
$ws.Range("A125").Value2 = 2050
$ws.Range("B125").Value2 = "2050_TM160_DBP_NoProject_06b"
$ws.Range("C125").Value2 = "RTP2025"
$ws.Range("D125").Value2 = "DBP"
$ws.Range("F125").Value2 = "Updated landuse/popsyn"
$ws.Range("G125").Value2 = "M:\urban_modeling\baus\PBA50Plus\PBA50Plus_NoProject_v7"
$ws.Range("H125").Value2 = "PBA50Plus_NoProject_v7"
$ws.Range("I125").Value2 = "current"
$ws.Range("K125").Value2 = "BlueprintNetworks_v18\net_2030_Baseline"
$ws.Range("L125").Value2 = "model3-b"
$ws.Range("M125").Value2 = "https://app.asana.com/0/1182463234225195/1207127941488239/f"
$ws.Range("N125").Value2 = 19.13
$ws.Range("O125").Value2 = "na"
$ws.Range("P125").Value2 = "na"
$ws.Range("Q125").Value2 = 0.87
$ws.Range("R125").Value2 = 0.78
$ws.Range("S125").Value2 = 83.3
$ws.Range("T125").Value2 = 0
$ws.Range("U125").Value2 = 72
$ws.Range("V125").Value2 = "2050 No Project"

# Demote the previous 2050 NoProject "current" row (still row 124)
$ws.Range("I124").Clear()
$ws.Range("V124").Clear()

# --- 2) 2035 Plan 07b -- insert after original row 117 (2035 Plan 06b) ---
$ws.Rows(118).Insert()
$ws.Range("A118").Value2 = 2035
$ws.Range("B118").Value2 = "2035_TM160_DBP_Plan_07b"
$ws.Range("C118").Value2 = "RTP2025"
$ws.Range("D118").Value2 = "DBP"
$ws.Range("F118").Value2 = "Update landuse/popsyn (fix too many persons)"
$ws.Range("G118").Value2 = "M:\urban_modeling\baus\PBA50Plus\PBA50Plus_DraftBlueprint\PBA50Plus_Draft_Blueprint_v6"
$ws.Range("H118").Value2 = "PBA50Plus_Draft_Blueprint_v6"
$ws.Range("I118").Value2 = "current"
$ws.Range("K118").Value2 = "BlueprintNetworks_v18\net_2035_Blueprint"
$ws.Range("L118").Value2 = "model2-b"
$ws.Range("M118").Value2 = "https://app.asana.com/0/1182463234225195/1207127280475842/f"
$ws.Range("N118").Value2 = 18.03
$ws.Range("O118").Value2 = "na"
$ws.Range("P118").Value2 = "na"
$ws.Range("Q118").Value2 = 0.87
$ws.Range("R118").Value2 = 0.78
$ws.Range("S118").Value2 = 83.3
$ws.Range("T118").Value2 = 0
$ws.Range("U118").Value2 = 72
$ws.Range("V118").Value2 = "2035 Plan"

# Demote the previous 2035 Plan "current" row (still row 117)
$ws.Range("I117").Clear()
$ws.Range("V117").Clear()

# --- 1) 2035 NoProject 06b -- insert after original row 110 (2035 NoProject 05b) ---
$ws.Rows(111).Insert()
$ws.Range("A111").Value2 = 2035
$ws.Range("B111").Value2 = "2035_TM160_DBP_NoProject_06b"
$ws.Range("C111").Value2 = "RTP2025"
$ws.Range("D111").Value2 = "DBP"
$ws.Range("F111").Value2 = "Updated landuse/popsyn"
$ws.Range("G111").Value2 = "M:\urban_modeling\baus\PBA50Plus\PBA50Plus_NoProject_v7"
$ws.Range("H111").Value2 = "PBA50Plus_NoProject_v7"
$ws.Range("I111").Value2 = "current"
$ws.Range("K111").Value2 = "BlueprintNetworks_v18\net_2030_Baseline"
$ws.Range("L111").Value2 = "model3-c"
$ws.Range("M111").Value2 = "https://app.asana.com/0/1182463234225195/1207127651751783/f"
$ws.Range("N111").Value2 = 17.579999999999998
$ws.Range("O111").Value2 = "na"
$ws.Range("P111").Value2 = "na"
$ws.Range("Q111").Value2 = 0.87
$ws.Range("R111").Value2 = 0.78
$ws.Range("S111").Value2 = 83.3
$ws.Range("T111").Value2 = 0
$ws.Range("U111").Value2 = 72
$ws.Range("V111").Value2 = "2035 No Project"

# Demote the previous 2035 NoProject "current" row (still row 110)
$ws.Range("I110").Clear()
$ws.Range("V110").Clear()
